# Apply the diff: insert two new rows at 287-288 (pushing the existing
# rows 287-393 down to 289-395), and populate the two new rows with the
# new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 287.
$ws.Range("A287:A288").EntireRow.Insert()

# --- New row 287 -----------------------------------------------------
$ws.Cells.Item(287, 1).Value = 7
$ws.Cells.Item(287, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(287, 3).Value = "Ñuble"
$ws.Cells.Item(287, 4).Value = 45119
$ws.Cells.Item(287, 5).Value = 16
$ws.Cells.Item(287, 6).Value = 100112017
$ws.Cells.Item(287, 7).Value = "Apio"
$ws.Cells.Item(287, 8).Value = "Americana (o)"
$ws.Cells.Item(287, 9).Value = "Primera"
$ws.Cells.Item(287, 10).Value = 150
$ws.Cells.Item(287, 11).Value = 7000
$ws.Cells.Item(287, 12).Value = 7000
$ws.Cells.Item(287, 13).Value = 7000
$ws.Cells.Item(287, 14).Value = "`$/docena de matas"
$ws.Cells.Item(287, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(287, 16).Value = 1167
$ws.Cells.Item(287, 17).Value = 6
$ws.Cells.Item(287, 18).Value = "Hortaliza"

# --- New row 288 -----------------------------------------------------
$ws.Cells.Item(288, 1).Value = 7
$ws.Cells.Item(288, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(288, 3).Value = "Ñuble"
$ws.Cells.Item(288, 4).Value = 45119
$ws.Cells.Item(288, 5).Value = 16
$ws.Cells.Item(288, 6).Value = 100112017
$ws.Cells.Item(288, 7).Value = "Apio"
$ws.Cells.Item(288, 8).Value = "Americana (o)"
$ws.Cells.Item(288, 9).Value = "Segunda"
$ws.Cells.Item(288, 10).Value = 100
$ws.Cells.Item(288, 11).Value = 6000
$ws.Cells.Item(288, 12).Value = 6000
$ws.Cells.Item(288, 13).Value = 6000
$ws.Cells.Item(288, 14).Value = "`$/docena de matas"
$ws.Cells.Item(288, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(288, 16).Value = 1000
$ws.Cells.Item(288, 17).Value = 6
$ws.Cells.Item(288, 18).Value = "Hortaliza"

# Make sure the D column on the two new rows keeps the date style used
# throughout the column (style index 2 in this workbook).
$ws.Range("D287:D288").NumberFormat = $ws.Range("D289").NumberFormat
